$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows after row 36 (before the old row 37), shifting all
# subsequent data down by 3 rows (old row 37 -> new row 40, ..., old row 89 -> new row 92).
$ws.Rows("37:39").Insert()

$ws.Range("A37").Value = 1
$ws.Range("B37").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C37").Value = 'Arica y Parinacota'
$ws.Range("D37").Value = 44966
$ws.Range("E37").Value = 15
$ws.Range("F37").Value = 'Fruta'
$ws.Range("G37").Value = 100103
$ws.Range("H37").Value = 'Frutos de hueso (carozo)'
$ws.Range("I37").Value = 100103006
$ws.Range("J37").Value = 'Nectarín'
$ws.Range("K37").Value = 'Artic Sprite'
$ws.Range("L37").Value = 'Primera'
$ws.Range("M37").Value = 300
$ws.Range("N37").Value = 24000
$ws.Range("O37").Value = 25000
$ws.Range("P37").Value = 24500
$ws.Range("Q37").Value = '$/bandeja 18 kilos granel'
$ws.Range("R37").Value = 'Región de O''Higgins'
$ws.Range("S37").Value = 1361
$ws.Range("T37").Value = 18

$ws.Range("A38").Value = 1
$ws.Range("B38").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C38").Value = 'Arica y Parinacota'
$ws.Range("D38").Value = 44966
$ws.Range("E38").Value = 15
$ws.Range("F38").Value = 'Fruta'
$ws.Range("G38").Value = 100103
$ws.Range("H38").Value = 'Frutos de hueso (carozo)'
$ws.Range("I38").Value = 100103006
$ws.Range("J38").Value = 'Nectarín'
$ws.Range("K38").Value = 'Super Queen'
$ws.Range("L38").Value = 'Segunda'
$ws.Range("M38").Value = 300
$ws.Range("N38").Value = 21000
$ws.Range("O38").Value = 22000
$ws.Range("P38").Value = 21500
$ws.Range("Q38").Value = '$/bandeja 18 kilos granel'
$ws.Range("R38").Value = 'Región de O''Higgins'
$ws.Range("S38").Value = 1194
$ws.Range("T38").Value = 18

$ws.Range("A39").Value = 1
$ws.Range("B39").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C39").Value = 'Arica y Parinacota'
$ws.Range("D39").Value = 44966
$ws.Range("E39").Value = 15
$ws.Range("F39").Value = 'Fruta'
$ws.Range("G39").Value = 100103
$ws.Range("H39").Value = 'Frutos de hueso (carozo)'
$ws.Range("I39").Value = 100103006
$ws.Range("J39").Value = 'Nectarín'
$ws.Range("K39").Value = 'Venus'
$ws.Range("L39").Value = 'Segunda'
$ws.Range("M39").Value = 270
$ws.Range("N39").Value = 21000
$ws.Range("O39").Value = 22000
$ws.Range("P39").Value = 21500
$ws.Range("Q39").Value = '$/bandeja 18 kilos granel'
$ws.Range("R39").Value = 'Región de O''Higgins'
$ws.Range("S39").Value = 1194
$ws.Range("T39").Value = 18
